$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newStatQuery = "MATCH (demo:demographic)`nWHERE demo.breed IN [`"Akita`"]`nMATCH (demo:demographic)-->(c:case)-->(s:study)-->(p:program)`nOPTIONAL MATCH (c)<-[*]-(samp:sample)`nOPTIONAL MATCH (c)<-[*]-(f:file)`nRETURN `n`tcount(DISTINCT(f)) as number_of_files, `n`tcount(DISTINCT(samp)) as number_of_sample, `n`tcount(DISTINCT(c)) as number_of_cases, `n`tcount(DISTINCT(s)) as number_of_study"

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery
